# Actualización última compra en mactrónica
# Updates purchase info for 4 components that were bought from "Mactrónica"
# (last purchased by "Danny"): rows 5, 6, 10 and 30 on sheet "Hoja1".
# The "estado compra" (F column) indicator flips from pending (red) to
# done (green) for each row now that quantity/price/buyer are filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$greenDone = 5287936   # RGB(0,176,80) - same "done" green used elsewhere in col F

# Item 4 (row 5): purchase store changes from "Centro" to "Mactrónica";
# quantity/unit price/buyer are now filled in.
$ws.Range("E5").Value = "Mactrónica"
$ws.Range("F5").Interior.Color = $greenDone
$ws.Range("G5").Value = 1500
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = "Danny"

# Item 5 (row 6): already "Mactrónica"; quantity/unit price/buyer filled in.
$ws.Range("F6").Interior.Color = $greenDone
$ws.Range("G6").Value = 2000
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = "Danny"

# Item 9 (row 10): already "Mactrónica"; quantity/unit price/buyer filled in.
$ws.Range("F10").Interior.Color = $greenDone
$ws.Range("G10").Value = 25000
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = "Danny"

# Item 29 (row 30): purchase store changes from "Centro" to "Mactrónica";
# quantity/unit price/buyer are now filled in.
$ws.Range("E30").Value = "Mactrónica"
$ws.Range("F30").Interior.Color = $greenDone
$ws.Range("G30").Value = 1500
$ws.Range("H30").Value = 1
$ws.Range("I30").Value = "Danny"
